# "Rename uninstall language keys"
#
# Rows 255-262 on the lang sheet hold the "Hub" uninstall-confirmation
# strings. Column A ("part") for every one of these rows becomes
# "uninstall" (it used to be a mix of confirm/remove/keep/deactivate/hub/msg),
# and several of the "var" (column B) keys get renamed:
#   remove_hub_data         (row255 B, was "confirm")      -> confirm_remove_hub_data
#   remove        -> (row257 A)                             -> uninstall
#   hub_data      (row257 B, unchanged text)
#   keep          -> (row258 A)                             -> uninstall
#                    (row258 B) keep                          -> keep_hub_data
#   deactivate/hub-> (row259 A/B)                             -> uninstall / deactivate_hub
#   hub/data      -> (row260 A/B)                             -> uninstall / hub_data
#   msg/removed_hub_data -> (row261 A/B)                      -> uninstall / msg_removed_hub_data
#   msg/kept_hub_data    -> (row262 A/B)                      -> uninstall / msg_kept_hub_data
#
# Row 256 (cancel/Cancel) is untouched - its text content does not change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 255: confirm / confirm_remove_hub_data / <unchanged de> / <unchanged en>
$ws.Range("A255").Value = "uninstall"
$ws.Range("B255").Value = "confirm_remove_hub_data"
$ws.Range("B255").WrapText = $true

# Row 257: remove / hub_data -> uninstall / remove_hub_data
$ws.Range("A257").Value = "uninstall"
$ws.Range("B257").Value = "remove_hub_data"

# Row 258: keep / hub_data -> uninstall / keep_hub_data
$ws.Range("A258").Value = "uninstall"
$ws.Range("B258").Value = "keep_hub_data"

# Row 259: deactivate / hub -> uninstall / deactivate_hub
$ws.Range("A259").Value = "uninstall"
$ws.Range("B259").Value = "deactivate_hub"
$ws.Range("B259").WrapText = $true

# Row 260: hub / data -> uninstall / hub_data
$ws.Range("A260").Value = "uninstall"
$ws.Range("B260").Value = "hub_data"

# Row 261: msg / removed_hub_data -> uninstall / msg_removed_hub_data
$ws.Range("A261").Value = "uninstall"
$ws.Range("B261").Value = "msg_removed_hub_data"

# Row 262: msg / kept_hub_data -> uninstall / msg_kept_hub_data
$ws.Range("A262").Value = "uninstall"
$ws.Range("B262").Value = "msg_kept_hub_data"

# Row heights for 257-262 nudge from 15.65 to 15.75 in the source file.
$ws.Rows.Item(257).RowHeight = 15.75
$ws.Rows.Item(258).RowHeight = 15.75
$ws.Rows.Item(259).RowHeight = 15.75
$ws.Rows.Item(260).RowHeight = 15.75
$ws.Rows.Item(261).RowHeight = 15.75
$ws.Rows.Item(262).RowHeight = 15.75

# Selection moved from B262 to A255 in the saved view state.
$ws.Range("A255").Select() | Out-Null
